$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 10
$ws.Range("B4").Value = 3

$ws.Range("A20").Value = "Auto number"
$ws.Range("B20").Value = "AutoNumber"
$ws.Range("C20").Value = "y"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5"

$ws.Range("B4").Select()
